# Automatische test-sync: 2025-08-28 20:34:50
# Appends a new "Retour status" log row to the Logs sheet and bumps the
# matching Dashboard category counter.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row goes right after the last existing data row (row 13 -> row 14).
$newRow = 14

$logs.Cells.Item($newRow, 1).Value = "Retour status"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($newRow, 6).Value = "2025-08-28 20:34:02"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Bump the Dashboard count for "Retour / Terugbetaling" from 12 to 13.
$dashboard.Cells.Item(2, 2).Value = 13

# Extend the conditional-formatting ranges (D/G/H/I/J) to cover the new row.
$oldLastRow = 13
$ccols = @("D", "G", "H", "I", "J")
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + $oldLastRow)
    $newRange = $logs.Range($col + "2:" + $col + $newRow)
    $fcs = $oldRange.FormatConditions
    if ($fcs.Count -gt 0) {
        $fcs.Item(1).ModifyAppliesToRange($newRange)
    }
}
